$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.380.15'
$ws.Range("E2").Value = '  -0.12%  '

# Row 3
$ws.Range("D3").Value = '1.846.87'
$ws.Range("E3").Value = '  -0.10%  '

# Row 4
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.27%  '

# Row 5
$ws.Range("D5").Value = '''240.66'
$ws.Range("E5").Value = '  +0.02%  '

# Row 6
$ws.Range("D6").Value = '''0.6287'
$ws.Range("E6").Value = '  +0.09%  '

# Row 7
$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '  +0.16%  '

# Row 8
$ws.Range("D8").Value = '''0.07491'
$ws.Range("E8").Value = '  -2.93%  '

# Row 9
$ws.Range("D9").Value = '''0.2895'
$ws.Range("E9").Value = '  -0.91%  '

# Row 10
$ws.Range("D10").Value = '''24.34'
$ws.Range("E10").Value = '  -2.93%  '

# Row 11
$ws.Range("D11").Value = '''0.07746'
$ws.Range("E11").Value = '  -0.01%  '

# Row 12
$ws.Range("D12").Value = '1.847.60'
$ws.Range("E12").Value = '  -0.22%  '

# Row 13
$ws.Range("D13").Value = '''5.018'
$ws.Range("E13").Value = '  -0.32%  '

# Row 14
$ws.Range("D14").Value = '''0.6800'
$ws.Range("E14").Value = '  -0.31%  '

# Row 15
$ws.Range("E15").Value = '  -4.75%  '

# Row 16
$ws.Range("D16").Value = '''82.98'
$ws.Range("E16").Value = '  -0.74%  '

# Row 17
$ws.Range("D17").Value = '2.114.94'
$ws.Range("E17").Value = '  -1.12%  '

# Row 18
$ws.Range("D18").Value = '''6.108'
$ws.Range("E18").Value = '  -1.32%  '

# Row 19
$ws.Range("D19").Value = '29.405.12'
$ws.Range("E19").Value = '  -0.15%  '

# Row 20
$ws.Range("D20").Value = '''229.02'
$ws.Range("E20").Value = '  +0.02%  '

# Row 21
$ws.Range("D21").Value = '''12.30'
$ws.Range("E21").Value = '  -0.73%  '

# Row 22
$ws.Range("D22").Value = '''1.002'
$ws.Range("E22").Value = '  +0.14%  '

# Row 23
$ws.Range("D23").Value = '''7.431'
$ws.Range("E23").Value = '  -0.32%  '

# Row 24
$ws.Range("E24").Value = '  +0.24%  '

# Row 25
$ws.Range("D25").Value = '''159.04'

# Row 26
$ws.Range("D26").Value = '''0.1384'
$ws.Range("E26").Value = '  +0.44%  '

# Row 27
$ws.Range("D27").Value = '''8.422'
$ws.Range("E27").Value = '  -0.08%  '

# Row 28
$ws.Range("D28").Value = '''17.58'
$ws.Range("E28").Value = '  -0.72%  '

# Row 29
$ws.Range("D29").Value = '''1.403'
$ws.Range("E29").Value = '  +4.10%  '

# Row 30
$ws.Range("D30").Value = '''1.477'
$ws.Range("E30").Value = '  +0.86%  '

# Row 31
$ws.Range("D31").Value = '''0.05689'
$ws.Range("E31").Value = '  +0.83%  '

# Row 32
$ws.Range("D32").Value = '''4.121'
$ws.Range("E32").Value = '  -0.15%  '

# Row 33
$ws.Range("D33").Value = '''4.044'
$ws.Range("E33").Value = '  -0.06%  '

# Row 34
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").Value = '''1.820'
$ws.Range("E34").Value = '  -1.35%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.150'
$ws.Range("E35").Value = '  -1.30%  '

# Row 36
$ws.Range("D36").Value = '''0.6941'
$ws.Range("E36").Value = '  -1.90%  '

# Row 37
$ws.Range("D37").Value = '''2.587'
$ws.Range("E37").Value = '  -0.27%  '

# Row 38
$ws.Range("D38").Value = '''2.837'
$ws.Range("E38").Value = '  +2.69%  '

# Row 39
$ws.Range("D39").Value = '1.249.78'
$ws.Range("E39").Value = '  +1.86%  '

# Row 40
$ws.Range("D40").Value = '''0.01820'
$ws.Range("E40").Value = '  +1.68%  '

# Row 41
$ws.Range("D41").Value = '''6.485'
$ws.Range("E41").Value = '  +0.50%  '

# Row 42
$ws.Range("D42").Value = '''0.9060'
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
$ws.Range("D43").Value = '''1.001'
$ws.Range("E43").Value = '  +0.00%  '

# Row 44
$ws.Range("D44").Value = '2.012.80'
$ws.Range("E44").Value = '  -1.52%  '

# Row 45
$ws.Range("D45").Value = '''101.29'
$ws.Range("E45").Value = '  -0.59%  '

# Row 46
$ws.Range("D46").Value = '''65.75'
$ws.Range("E46").Value = '  -0.56%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '''7.074'
$ws.Range("E47").Value = '  -1.68%  '

# Row 48
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '''0.1159'
$ws.Range("E48").Value = '  +0.10%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '''0.00000000116'
$ws.Range("E49").Value = '  -5.09%  '

# Row 50
$ws.Range("D50").Value = '''9.002'
$ws.Range("E50").Value = '  -0.06%  '

# Row 51
$ws.Range("D51").Value = '''0.3940'
$ws.Range("E51").Value = '  -2.07%  '
